$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-order match rows within same-day groups (home/away + odds + url moved between rows) ---
# Row 44
$ws.Range("F44").Value = 'Pennarossa'
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = 'Domagnano'
$ws.Range("I44").Value = 2
$ws.Range("J44").Value = 2.35
$ws.Range("K44").Value = '27/10/2023 03:12'
$ws.Range("L44").Value = 3.25
$ws.Range("M44").Value = '28/10/2023 14:46'
$ws.Range("N44").Value = 3.1
$ws.Range("O44").Value = '27/10/2023 03:12'
$ws.Range("P44").Value = 3.66
$ws.Range("Q44").Value = '28/10/2023 14:33'
$ws.Range("R44").Value = 2.56
$ws.Range("S44").Value = '27/10/2023 03:12'
$ws.Range("T44").Value = 1.9
$ws.Range("U44").Value = '28/10/2023 14:46'
$ws.Range("V44").Value = 'https://www.betexplorer.com/football/san-marino/campionato-sammarinese/ss-pennarossa-sp-domagnano/UevtEPAf/'

# Row 45
$ws.Range("F45").Value = 'San Giovanni'
$ws.Range("G45").Value = 2
$ws.Range("H45").Value = 'Fiorentino'
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 3.06
$ws.Range("K45").Value = '27/10/2023 03:12'
$ws.Range("L45").Value = 3.24
$ws.Range("M45").Value = '28/10/2023 14:33'
$ws.Range("N45").Value = 3.21
$ws.Range("O45").Value = '27/10/2023 03:12'
$ws.Range("P45").Value = 3.24
$ws.Range("Q45").Value = '28/10/2023 14:33'
$ws.Range("R45").Value = 1.96
$ws.Range("S45").Value = '27/10/2023 03:12'
$ws.Range("T45").Value = 2.05
$ws.Range("U45").Value = '28/10/2023 14:33'
$ws.Range("V45").Value = 'https://www.betexplorer.com/football/san-marino/campionato-sammarinese/san-giovanni-fiorentino/xba1L7IK/'

# Row 46
$ws.Range("F46").Value = 'Virtus'
$ws.Range("G46").Value = 2
$ws.Range("H46").Value = 'La Fiorita'
$ws.Range("I46").Value = 1
$ws.Range("J46").Value = 2.64
$ws.Range("K46").Value = '27/10/2023 03:12'
$ws.Range("L46").Value = 2.82
$ws.Range("M46").Value = '28/10/2023 14:33'
$ws.Range("N46").Value = 2.75
$ws.Range("O46").Value = '27/10/2023 03:12'
$ws.Range("P46").Value = 2.69
$ws.Range("Q46").Value = '28/10/2023 14:33'
$ws.Range("R46").Value = 2.47
$ws.Range("S46").Value = '27/10/2023 03:12'
$ws.Range("T46").Value = 2.64
$ws.Range("U46").Value = '28/10/2023 14:33'
$ws.Range("V46").Value = 'https://www.betexplorer.com/football/san-marino/campionato-sammarinese/virtus-la-fiorita/Glb5KRXQ/'

# Row 48
$ws.Range("F48").Value = 'Tre Fiori'
$ws.Range("G48").Value = 5
$ws.Range("H48").Value = 'Cosmos'
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 3.11
$ws.Range("K48").Value = '28/10/2023 04:13'
$ws.Range("L48").Value = 2.84
$ws.Range("M48").Value = '29/10/2023 09:05'
$ws.Range("N48").Value = 2.97
$ws.Range("O48").Value = '28/10/2023 04:13'
$ws.Range("P48").Value = 3.25
$ws.Range("Q48").Value = '29/10/2023 13:01'
$ws.Range("R48").Value = 2.05
$ws.Range("S48").Value = '28/10/2023 04:13'
$ws.Range("T48").Value = 2.2
$ws.Range("U48").Value = '29/10/2023 09:05'
$ws.Range("V48").Value = 'https://www.betexplorer.com/football/san-marino/campionato-sammarinese/tre-fiori-sp-cosmos/8OfcMm3E/'

# Row 49
$ws.Range("F49").Value = 'Cailungo'
$ws.Range("G49").Value = 4
$ws.Range("H49").Value = 'Faetano'
$ws.Range("I49").Value = 1
$ws.Range("J49").Value = 2.69
$ws.Range("K49").Value = '28/10/2023 04:13'
$ws.Range("L49").Value = 3.01
$ws.Range("M49").Value = '29/10/2023 14:54'
$ws.Range("N49").Value = 3.3
$ws.Range("O49").Value = '28/10/2023 04:13'
$ws.Range("P49").Value = 3.7
$ws.Range("Q49").Value = '29/10/2023 14:54'
$ws.Range("R49").Value = 2.12
$ws.Range("S49").Value = '28/10/2023 04:13'
$ws.Range("T49").Value = 1.99
$ws.Range("U49").Value = '29/10/2023 14:54'
$ws.Range("V49").Value = 'https://www.betexplorer.com/football/san-marino/campionato-sammarinese/cailungo-sc-faetano/MVzYFotr/'

# Row 50
$ws.Range("F50").Value = 'Tre Penne'
$ws.Range("G50").Value = 3
$ws.Range("H50").Value = 'Murata'
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 1.45
$ws.Range("K50").Value = '28/10/2023 04:13'
$ws.Range("L50").Value = 1.85
$ws.Range("M50").Value = '29/10/2023 14:24'
$ws.Range("N50").Value = 4.07
$ws.Range("O50").Value = '28/10/2023 04:13'
$ws.Range("P50").Value = 3.81
$ws.Range("Q50").Value = '29/10/2023 14:32'
$ws.Range("R50").Value = 4.66
$ws.Range("S50").Value = '28/10/2023 04:13'
$ws.Range("T50").Value = 3.26
$ws.Range("U50").Value = '29/10/2023 14:24'
$ws.Range("V50").Value = 'https://www.betexplorer.com/football/san-marino/campionato-sammarinese/tre-penne-ss-murata/zZvxF5el/'

# Row 70
$ws.Range("F70").Value = 'Tre Penne'
$ws.Range("G70").Value = 2
$ws.Range("H70").Value = 'Libertas'
$ws.Range("I70").Value = 1
$ws.Range("J70").Value = 1.3
$ws.Range("K70").Value = '02/12/2023 04:12'
$ws.Range("L70").Value = 1.29
$ws.Range("M70").Value = '02/12/2023 14:35'
$ws.Range("N70").Value = 4.86
$ws.Range("O70").Value = '02/12/2023 04:12'
$ws.Range("P70").Value = 5.03
$ws.Range("Q70").Value = '02/12/2023 14:35'
$ws.Range("R70").Value = 7.16
$ws.Range("S70").Value = '02/12/2023 04:12'
$ws.Range("T70").Value = 7.76
$ws.Range("U70").Value = '02/12/2023 14:35'
$ws.Range("V70").Value = 'https://www.betexplorer.com/football/san-marino/campionato-sammarinese/tre-penne-ac-libertas/rcgxvDUG/'

# Row 71
$ws.Range("F71").Value = 'Virtus'
$ws.Range("G71").Value = 2
$ws.Range("H71").Value = 'Murata'
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 1.78
$ws.Range("K71").Value = '02/12/2023 04:12'
$ws.Range("L71").Value = 1.73
$ws.Range("M71").Value = '02/12/2023 14:44'
$ws.Range("N71").Value = 3.55
$ws.Range("O71").Value = '02/12/2023 04:12'
$ws.Range("P71").Value = 3.47
$ws.Range("Q71").Value = '02/12/2023 14:44'
$ws.Range("R71").Value = 3.68
$ws.Range("S71").Value = '02/12/2023 04:12'
$ws.Range("T71").Value = 4.17
$ws.Range("U71").Value = '02/12/2023 13:52'
$ws.Range("V71").Value = 'https://www.betexplorer.com/football/san-marino/campionato-sammarinese/virtus-ss-murata/x0jUuZa4/'

# Row 80
$ws.Range("F80").Value = 'La Fiorita'
$ws.Range("G80").Value = 3
$ws.Range("H80").Value = 'San Marino Academy U22'
$ws.Range("I80").Value = 1
$ws.Range("J80").Value = 1.18
$ws.Range("K80").Value = '09/12/2023 12:12'
$ws.Range("L80").Value = 1.15
$ws.Range("M80").Value = '09/12/2023 13:23'
$ws.Range("N80").Value = 6.1
$ws.Range("O80").Value = '09/12/2023 12:12'
$ws.Range("P80").Value = 6.7
$ws.Range("Q80").Value = '09/12/2023 18:03'
$ws.Range("R80").Value = 9.96
$ws.Range("S80").Value = '09/12/2023 12:12'
$ws.Range("T80").Value = 11.6
$ws.Range("U80").Value = '09/12/2023 18:03'
$ws.Range("V80").Value = 'https://www.betexplorer.com/football/san-marino/campionato-sammarinese/la-fiorita-san-marino-academy/4SV8iWbG/'

# Row 81
$ws.Range("F81").Value = 'San Giovanni'
$ws.Range("G81").Value = 5
$ws.Range("H81").Value = 'Cailungo'
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 2.37
$ws.Range("K81").Value = '09/12/2023 07:42'
$ws.Range("L81").Value = 2.28
$ws.Range("M81").Value = '09/12/2023 16:20'
$ws.Range("N81").Value = 3.19
$ws.Range("O81").Value = '09/12/2023 07:42'
$ws.Range("P81").Value = 3.24
$ws.Range("Q81").Value = '09/12/2023 17:58'
$ws.Range("R81").Value = 2.65
$ws.Range("S81").Value = '09/12/2023 07:42'
$ws.Range("T81").Value = 2.79
$ws.Range("U81").Value = '09/12/2023 17:58'
$ws.Range("V81").Value = 'https://www.betexplorer.com/football/san-marino/campionato-sammarinese/san-giovanni-cailungo/xtUGkATS/'

# Row 85
$ws.Range("F85").Value = 'Pennarossa'
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 'San Giovanni'
$ws.Range("I85").Value = 3
$ws.Range("J85").Value = 3.28
$ws.Range("K85").Value = '16/12/2023 04:13'
$ws.Range("L85").Value = 6.29
$ws.Range("M85").Value = '16/12/2023 14:59'
$ws.Range("N85").Value = 3.31
$ws.Range("O85").Value = '16/12/2023 04:13'
$ws.Range("P85").Value = 4.32
$ws.Range("Q85").Value = '16/12/2023 14:59'
$ws.Range("R85").Value = 1.97
$ws.Range("S85").Value = '16/12/2023 04:13'
$ws.Range("T85").Value = 1.4
$ws.Range("U85").Value = '16/12/2023 14:59'
$ws.Range("V85").Value = 'https://www.betexplorer.com/football/san-marino/campionato-sammarinese/ss-pennarossa-san-giovanni/OW2Z7iqc/'

# Row 86
$ws.Range("F86").Value = 'Juvenes/Dogana'
$ws.Range("G86").Value = 3
$ws.Range("H86").Value = 'Faetano'
$ws.Range("I86").Value = 1
$ws.Range("J86").Value = 2.08
$ws.Range("K86").Value = '16/12/2023 04:13'
$ws.Range("L86").Value = 2.27
$ws.Range("M86").Value = '16/12/2023 13:00'
$ws.Range("N86").Value = 3.5
$ws.Range("O86").Value = '16/12/2023 04:13'
$ws.Range("P86").Value = 3.6
$ws.Range("Q86").Value = '16/12/2023 13:00'
$ws.Range("R86").Value = 2.88
$ws.Range("S86").Value = '16/12/2023 04:13'
$ws.Range("T86").Value = 2.58
$ws.Range("U86").Value = '16/12/2023 13:00'
$ws.Range("V86").Value = 'https://www.betexplorer.com/football/san-marino/campionato-sammarinese/juvenes-dogana-sc-faetano/2N1w7Bb3/'

# Row 87
$ws.Range("F87").Value = 'Cosmos'
$ws.Range("G87").Value = 2
$ws.Range("H87").Value = 'Murata'
$ws.Range("I87").Value = 1
$ws.Range("J87").Value = 1.64
$ws.Range("K87").Value = '16/12/2023 04:13'
$ws.Range("L87").Value = 1.54
$ws.Range("M87").Value = '16/12/2023 14:59'
$ws.Range("N87").Value = 3.77
$ws.Range("O87").Value = '16/12/2023 04:13'
$ws.Range("P87").Value = 3.9
$ws.Range("Q87").Value = '16/12/2023 14:59'
$ws.Range("R87").Value = 4.17
$ws.Range("S87").Value = '16/12/2023 04:13'
$ws.Range("T87").Value = 4.95
$ws.Range("U87").Value = '16/12/2023 14:59'
$ws.Range("V87").Value = 'https://www.betexplorer.com/football/san-marino/campionato-sammarinese/sp-cosmos-ss-murata/zTXEmNcr/'

# Row 89
$ws.Range("F89").Value = 'Fiorentino'
$ws.Range("G89").Value = 1
$ws.Range("H89").Value = 'Folgore'
$ws.Range("I89").Value = 2
$ws.Range("J89").Value = 2.86
$ws.Range("K89").Value = '17/12/2023 04:12'
$ws.Range("L89").Value = 2.92
$ws.Range("M89").Value = '17/12/2023 14:46'
$ws.Range("N89").Value = 3.24
$ws.Range("O89").Value = '17/12/2023 04:12'
$ws.Range("P89").Value = 3.22
$ws.Range("Q89").Value = '17/12/2023 14:46'
$ws.Range("R89").Value = 2.2
$ws.Range("S89").Value = '17/12/2023 04:12'
$ws.Range("T89").Value = 2.21
$ws.Range("U89").Value = '17/12/2023 14:46'
$ws.Range("V89").Value = 'https://www.betexplorer.com/football/san-marino/campionato-sammarinese/fiorentino-folgore/AoQRpLs1/'

# Row 90
$ws.Range("F90").Value = 'Tre Penne'
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 'La Fiorita'
$ws.Range("I90").Value = 2
$ws.Range("J90").Value = 2.44
$ws.Range("K90").Value = '17/12/2023 04:12'
$ws.Range("L90").Value = 3
$ws.Range("M90").Value = '17/12/2023 14:49'
$ws.Range("N90").Value = 2.87
$ws.Range("O90").Value = '17/12/2023 04:12'
$ws.Range("P90").Value = 2.36
$ws.Range("Q90").Value = '17/12/2023 14:49'
$ws.Range("R90").Value = 2.81
$ws.Range("S90").Value = '17/12/2023 04:12'
$ws.Range("T90").Value = 2.87
$ws.Range("U90").Value = '17/12/2023 14:49'
$ws.Range("V90").Value = 'https://www.betexplorer.com/football/san-marino/campionato-sammarinese/tre-penne-la-fiorita/EJWInsDl/'

# Row 91
$ws.Range("F91").Value = 'Virtus'
$ws.Range("G91").Value = 3
$ws.Range("H91").Value = 'Cailungo'
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 1.14
$ws.Range("K91").Value = '17/12/2023 04:12'
$ws.Range("L91").Value = 1.13
$ws.Range("M91").Value = '17/12/2023 14:02'
$ws.Range("N91").Value = 6.85
$ws.Range("O91").Value = '17/12/2023 04:12'
$ws.Range("P91").Value = 7.4
$ws.Range("Q91").Value = '17/12/2023 14:02'
$ws.Range("R91").Value = 11.07
$ws.Range("S91").Value = '17/12/2023 04:12'
$ws.Range("T91").Value = 12.37
$ws.Range("U91").Value = '17/12/2023 14:02'
$ws.Range("V91").Value = 'https://www.betexplorer.com/football/san-marino/campionato-sammarinese/virtus-cailungo/pnir6VD9/'

# --- Append 3 new match rows (96:98), copying formatting from the last existing row (95) ---
$ws.Range("A95:V95").Copy()
$ws.Range("A96:V98").PasteSpecial(-4122)

# Row 96
$ws.Range("A96").Value = 95
$ws.Range("B96").Value = 'san-marino'
$ws.Range("C96").Value = 'campionato-sammarinese'
$ws.Range("D96").Value = '2023-2024'
$ws.Range("E96").Value = 45298.625
$ws.Range("F96").Value = 'Folgore'
$ws.Range("G96").Value = 8
$ws.Range("H96").Value = 'Pennarossa'
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 1.44
$ws.Range("K96").Value = '07/01/2024 04:12'
$ws.Range("L96").Value = 1.37
$ws.Range("M96").Value = '07/01/2024 14:28'
$ws.Range("N96").Value = 4.19
$ws.Range("O96").Value = '07/01/2024 04:12'
$ws.Range("P96").Value = 4.31
$ws.Range("Q96").Value = '07/01/2024 14:30'
$ws.Range("R96").Value = 5.5
$ws.Range("S96").Value = '07/01/2024 04:12'
$ws.Range("T96").Value = 6.97
$ws.Range("U96").Value = '07/01/2024 14:30'
$ws.Range("V96").Value = 'https://www.betexplorer.com/football/san-marino/campionato-sammarinese/folgore-ss-pennarossa/b76M624K/'

# Row 97
$ws.Range("A97").Value = 96
$ws.Range("B97").Value = 'san-marino'
$ws.Range("C97").Value = 'campionato-sammarinese'
$ws.Range("D97").Value = '2023-2024'
$ws.Range("E97").Value = 45298.625
$ws.Range("F97").Value = 'Libertas'
$ws.Range("G97").Value = 1
$ws.Range("H97").Value = 'Cosmos'
$ws.Range("I97").Value = 3
$ws.Range("J97").Value = 6.33
$ws.Range("K97").Value = '07/01/2024 04:12'
$ws.Range("L97").Value = 7.27
$ws.Range("M97").Value = '07/01/2024 14:04'
$ws.Range("N97").Value = 4.48
$ws.Range("O97").Value = '07/01/2024 04:12'
$ws.Range("P97").Value = 5.07
$ws.Range("Q97").Value = '07/01/2024 14:04'
$ws.Range("R97").Value = 1.37
$ws.Range("S97").Value = '07/01/2024 04:12'
$ws.Range("T97").Value = 1.3
$ws.Range("U97").Value = '07/01/2024 13:39'
$ws.Range("V97").Value = 'https://www.betexplorer.com/football/san-marino/campionato-sammarinese/ac-libertas-sp-cosmos/KG7I7rkE/'

# Row 98
$ws.Range("A98").Value = 97
$ws.Range("B98").Value = 'san-marino'
$ws.Range("C98").Value = 'campionato-sammarinese'
$ws.Range("D98").Value = '2023-2024'
$ws.Range("E98").Value = 45298.625
$ws.Range("F98").Value = 'San Giovanni'
$ws.Range("G98").Value = 1
$ws.Range("H98").Value = 'Virtus'
$ws.Range("I98").Value = 2
$ws.Range("J98").Value = 10.07
$ws.Range("K98").Value = '07/01/2024 04:12'
$ws.Range("L98").Value = 6.67
$ws.Range("M98").Value = '07/01/2024 13:03'
$ws.Range("N98").Value = 5.58
$ws.Range("O98").Value = '07/01/2024 04:12'
$ws.Range("P98").Value = 4.49
$ws.Range("Q98").Value = '07/01/2024 13:03'
$ws.Range("R98").Value = 1.2
$ws.Range("S98").Value = '07/01/2024 04:12'
$ws.Range("T98").Value = 1.36
$ws.Range("U98").Value = '07/01/2024 13:03'
$ws.Range("V98").Value = 'https://www.betexplorer.com/football/san-marino/campionato-sammarinese/san-giovanni-virtus/48PVquc7/'

Write-Output "done"